$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 631 (shifts existing rows 631..684 down to 632..685)
$ws.Rows("631:631").Insert()

# Fill in the new row 631 with data (same constants as neighboring rows, new specifics)
$ws.Range("A631").Value = 10
$ws.Range("B631").Value = "Vega Modelo de Temuco"
$ws.Range("C631").Value = "La Araucanía"
$ws.Range("D631").Value = 44769
$ws.Range("E631").Value = 9
$ws.Range("F631").Value = "Fruta"
$ws.Range("G631").Value = 100108
$ws.Range("H631").Value = "Tropicales y subtropicales"
$ws.Range("I631").Value = 100108006
$ws.Range("J631").Value = "Plátano"
$ws.Range("K631").Value = "Sin especificar"
$ws.Range("L631").Value = "Pintón"
$ws.Range("M631").Value = 480
$ws.Range("N631").Value = 34000
$ws.Range("O631").Value = 34000
$ws.Range("P631").Value = 34000
$ws.Range("Q631").Value = "$/caja 20 kilos"
$ws.Range("R631").Value = "Ecuador"
$ws.Range("S631").Value = 1700
$ws.Range("T631").Value = 20
